# chore: update Sheets via scheduled runner
# Refresh market-board derived price/profit columns (H-N) for a batch of
# leve rows across the crafting-job sheets, as produced by the scheduled
# price-update runner.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 180
$ws.Range("I31").Value = 180
$ws.Range("K31").Value = 540
$ws.Range("M31").Value = -310
$ws.Range("H116").Value = 5475
$ws.Range("I116").Value = 5962
$ws.Range("J116").Value = 4988
$ws.Range("K116").Value = 5962
$ws.Range("L116").Value = 4988
$ws.Range("M116").Value = -2520
$ws.Range("N116").Value = -11872

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 487.64706
$ws.Range("I80").Value = 438
$ws.Range("J80").Value = 531.7778
$ws.Range("K80").Value = 438
$ws.Range("L80").Value = 531.7778
$ws.Range("M80").Value = 560
$ws.Range("N80").Value = -2527.7778
$ws.Range("H83").Value = 487.64706
$ws.Range("I83").Value = 438
$ws.Range("J83").Value = 531.7778
$ws.Range("K83").Value = 2190
$ws.Range("L83").Value = 2658.889
$ws.Range("M83").Value = 2802
$ws.Range("N83").Value = -12642.889
$ws.Range("H86").Value = 2392.7693
$ws.Range("I86").Value = 1976.5
$ws.Range("K86").Value = 1976.5
$ws.Range("M86").Value = -853.5
$ws.Range("H89").Value = 2392.7693
$ws.Range("I89").Value = 1976.5
$ws.Range("K89").Value = 9882.5
$ws.Range("M89").Value = -4266.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H43").Value = 5380
$ws.Range("J43").Value = 5380
$ws.Range("L43").Value = 5380
$ws.Range("N43").Value = -5748
$ws.Range("H68").Value = 19860
$ws.Range("J68").Value = 21450
$ws.Range("L68").Value = 21450
$ws.Range("N68").Value = -22948
$ws.Range("H71").Value = 19860
$ws.Range("J71").Value = 21450
$ws.Range("L71").Value = 64350
$ws.Range("N71").Value = -71838
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").ClearContents()
$ws.Range("N80").Value = 0
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").ClearContents()
$ws.Range("N83").Value = 0
$ws.Range("H87").Value = 22332.666
$ws.Range("J87").Value = 22332.666
$ws.Range("L87").Value = 22332.666
$ws.Range("N87").Value = -24704.666
$ws.Range("H90").Value = 22332.666
$ws.Range("J90").Value = 22332.666
$ws.Range("L90").Value = 66997.99800000001
$ws.Range("N90").Value = -78853.99800000001
$ws.Range("H93").Value = 4209.727
$ws.Range("I93").Value = 4209.727
$ws.Range("K93").Value = 4209.727
$ws.Range("M93").Value = -2337.727
$ws.Range("H95").Value = 333356540
$ws.Range("J95").Value = 333356540
$ws.Range("L95").Value = 333356540
$ws.Range("N95").Value = -333362032
$ws.Range("H96").Value = 200021490
$ws.Range("J96").Value = 200021490
$ws.Range("L96").Value = 200021490
$ws.Range("N96").Value = -200026982
$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").ClearContents()
$ws.Range("N97").Value = 0
$ws.Range("H101").Value = 5380
$ws.Range("J101").Value = 5380
$ws.Range("L101").Value = 5380
$ws.Range("N101").Value = -11870
$ws.Range("H102").Value = 29550
$ws.Range("J102").Value = 29550
$ws.Range("L102").Value = 29550
$ws.Range("N102").Value = -34418
$ws.Range("H103").Value = 6124.8
$ws.Range("I103").Value = 7081
$ws.Range("J103").Value = 2300
$ws.Range("K103").Value = 7081
$ws.Range("L103").Value = 2300
$ws.Range("M103").Value = -5909
$ws.Range("N103").Value = -4644
$ws.Range("H104").Value = 24900
$ws.Range("J104").Value = 24900
$ws.Range("L104").Value = 24900
$ws.Range("N104").Value = -30142

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 774.03845
$ws.Range("I4").Value = 148.77777
$ws.Range("J4").Value = 1105.0588
$ws.Range("K4").Value = 446.33331
$ws.Range("L4").Value = 3315.1764
$ws.Range("M4").Value = -334.33331
$ws.Range("N4").Value = -3539.1764
$ws.Range("H6").Value = 109
$ws.Range("I6").Value = 54.444443
$ws.Range("J6").Value = 600
$ws.Range("K6").Value = 163.333329
$ws.Range("L6").Value = 1800
$ws.Range("M6").Value = -50.33332899999999
$ws.Range("N6").Value = -2026
$ws.Range("H44").Value = 30303488
$ws.Range("I44").Value = 255.83333
$ws.Range("J44").Value = 66667370
$ws.Range("K44").Value = 767.49999
$ws.Range("L44").Value = 200002110
$ws.Range("M44").Value = -369.49999
$ws.Range("N44").Value = -200002906
$ws.Range("H54").Value = 4000
$ws.Range("J54").Value = 4000
$ws.Range("L54").Value = 12000
$ws.Range("N54").Value = -13118
$ws.Range("H63").Value = 2685.6667
$ws.Range("I63").Value = 1675
$ws.Range("J63").Value = 4707
$ws.Range("K63").Value = 5025
$ws.Range("L63").Value = 14121
$ws.Range("M63").Value = -4276
$ws.Range("N63").Value = -15619
$ws.Range("H64").Value = 1260
$ws.Range("I64").Value = 850
$ws.Range("K64").Value = 2550
$ws.Range("M64").Value = -2280
$ws.Range("H66").Value = 2685.6667
$ws.Range("I66").Value = 1675
$ws.Range("J66").Value = 4707
$ws.Range("K66").Value = 15075
$ws.Range("L66").Value = 42363
$ws.Range("M66").Value = -11331
$ws.Range("N66").Value = -49851
$ws.Range("H67").Value = 1260
$ws.Range("I67").Value = 850
$ws.Range("K67").Value = 2550
$ws.Range("M67").Value = -1614
$ws.Range("H120").Value = 11616.5
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 11616.5
$ws.Range("K120").Value = 0
$ws.Range("L120").ClearContents()
$ws.Range("M120").Value = 34849.5
$ws.Range("N120").Value = -44525.5
$ws.Range("H131").Value = 641.4857
$ws.Range("I131").Value = 337.94116
$ws.Range("J131").Value = 928.1667
$ws.Range("K131").Value = 1013.82348
$ws.Range("L131").Value = 2784.5001
$ws.Range("M131").Value = 4026.17652
$ws.Range("N131").Value = -12864.5001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 62275.176
$ws.Range("I122").Value = 114066.445
$ws.Range("J122").Value = 4010
$ws.Range("K122").Value = 342199.335
$ws.Range("L122").Value = 12030
$ws.Range("M122").Value = -339749.335
$ws.Range("N122").Value = -16930

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 37110.332
$ws.Range("J80").Value = 37110.332
$ws.Range("L80").Value = 37110.332
$ws.Range("N80").Value = -39106.332
$ws.Range("H83").Value = 37110.332
$ws.Range("J83").Value = 37110.332
$ws.Range("L83").Value = 111330.996
$ws.Range("N83").Value = -121314.996
